# Add the new "Class Information" worksheet as the LAST sheet in the workbook
# (Worksheets.Add() with no args inserts before the last sheet, so we pass an
# explicit After: = the current last worksheet).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Class Information"

# Header row
$ws.Range("A1").Value = "Class Name"
$ws.Range("B1").Value = "Time Spend"
$ws.Range("C1").Value = "Class ID"
$ws.Range("D1").Value = "Professor"
$ws.Range("E1").Value = "Room"
$ws.Range("F1").Value = "Max Occupancy"
$ws.Range("G1").Value = "Current Occupancy"

# Row 2 - English
$ws.Range("A2").Value = "English"
$ws.Range("B2").Value = "8:00-9:00"
$ws.Range("C2").Value = 12345
$ws.Range("D2").Value = "Jane Doe"
$ws.Range("E2").Value = "A-123"
$ws.Range("F2").Value = 30

# Row 3 - Geometry
$ws.Range("A3").Value = "Geometry"
$ws.Range("B3").Value = "9:00-10:00"
$ws.Range("C3").Value = 12346
$ws.Range("D3").Value = "Kevin Smith"
$ws.Range("E3").Value = "A-231"
$ws.Range("F3").Value = 30

# Row 4 - Physics
$ws.Range("A4").Value = "Physics"
$ws.Range("B4").Value = "10:00-11:00"
$ws.Range("C4").Value = 12347
$ws.Range("D4").Value = "Jennifer Jones"
$ws.Range("E4").Value = "B-102"
$ws.Range("F4").Value = 30

# Row 5 - Chemistry (Time Spend cell carries an explicit h:mm number format)
$ws.Range("A5").Value = "Chemistry"
$ws.Range("B5").Value = "11:00-12:00"
$ws.Range("C5").Value = 12348
$ws.Range("D5").Value = "David Smith"
$ws.Range("E5").Value = "B-204"
$ws.Range("F5").Value = 25
$ws.Range("B5").NumberFormat = "h:mm"

# Row 6 - Biology
$ws.Range("A6").Value = "Biology"
$ws.Range("B6").Value = "12:00-13:00"
$ws.Range("C6").Value = 12349
$ws.Range("D6").Value = "Eric Huang"
$ws.Range("E6").Value = "B-123"
$ws.Range("F6").Value = 25

# Row 7 - Philosphy
$ws.Range("A7").Value = "Philosphy"
$ws.Range("B7").Value = "13:00-14:00"
$ws.Range("C7").Value = 12350
$ws.Range("D7").Value = "Vincent Ku"
$ws.Range("E7").Value = "A-123"
$ws.Range("F7").Value = 30

# Row 8 - Linear Algebra
$ws.Range("A8").Value = "Linear Algebra"
$ws.Range("B8").Value = "14:00-15:00"
$ws.Range("C8").Value = 12351
$ws.Range("D8").Value = "Joshua Hu"
$ws.Range("E8").Value = "A-231"
$ws.Range("F8").Value = 30

# Row 9 - Discrete Math
$ws.Range("A9").Value = "Discrete Math"
$ws.Range("B9").Value = "15:00-16:00"
$ws.Range("C9").Value = 12352
$ws.Range("D9").Value = "Linda J"
$ws.Range("E9").Value = "B-102"
$ws.Range("F9").Value = 30

# Row 10 - Calculus
$ws.Range("A10").Value = "Calculus"
$ws.Range("B10").Value = "16:00-17:00"
$ws.Range("C10").Value = 12353
$ws.Range("D10").Value = "Henry H"
$ws.Range("E10").Value = "B-204"
$ws.Range("F10").Value = 25

# Row 11 - Python
$ws.Range("A11").Value = "Python"
$ws.Range("B11").Value = "17:00-18:00"
$ws.Range("C11").Value = 12354
$ws.Range("D11").Value = "Julie J"
$ws.Range("E11").Value = "B-123"
$ws.Range("F11").Value = 25

# Row 12 - Java
$ws.Range("A12").Value = "Java"
$ws.Range("B12").Value = "18:00-19:00"
$ws.Range("C12").Value = 12355
$ws.Range("D12").Value = "Jason L"
$ws.Range("E12").Value = "A-210"
$ws.Range("F12").Value = 30

# Match the author's last selection on the new sheet
$ws.Range("H21").Select() | Out-Null
